$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: bump the date in A1 by one day
$ws.Range("A1").Value = 45309

# Step 2: update the price list values in column D
$ws.Range("D33").Value = 1165.8
$ws.Range("D34").Value = 1966.2
$ws.Range("D35").Value = 2164.56
